# Add new "2m" wind variable legend rows to the "Taula llegendes" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Insert "VAR_DVM_2_m_graus" row right before "VAR_GN_cm" (old row 16) ---
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "VAR_DVM_2_m_graus       "
$ws.Range("B16").Value = " Direcció vent 2m"

# --- 2) Insert "VAR_VVM_2_m_km_h" row right after "VAR_VVM_6_m_km_h" (now row 27) ---
$ws.Rows.Item(28).Insert()
$ws.Range("A28").Value = "VAR_VVM_2_m_km_h   "
$ws.Range("B28").Value = " Velocitat Mitjana del Vent (Km/h)"

# --- 3) Append "VAR_VVX_2_m_km_h" row at the very end (new row 31) ---
$ws.Range("A31").Value = "VAR_VVX_2_m_km_h     "
$ws.Range("B31").Value = " Ratxa Màxima del Vent (Km/h)"

# --- 4) Update the view: scrolled down a bit, selection on B31 ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("B31").Select()
